$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The B column holds price strings that look numeric (e.g. "11,999"); force
# them to stay plain text (matching the source workbook, which stores every
# price as a shared string) instead of letting Excel auto-convert them to
# numbers with a thousands-separator number format.
$priceRange = $ws.Range("B2:B11")
$priceRange.NumberFormat = "@"

# Updated product listing (rows 2-6), replacing the previous monitor lineup
$ws.Range("A2").Value = 'Lenovo Monitor Legion R27fc-30 Gaming Curved Monitor, 27" FHD VA Display, Up to 280Hz Refresh Rate, 0.5ms (MPRT) Response Time, 1500R Curvature, FreeSync, Adaptive Sync & G-Sync Compatible, Black'
$ws.Range("B2").Value = "11,999"

$ws.Range("A3").Value = "Samsung 27-Inch G55C Odyssey QHD 2K Curved Gaming Monitor, HRDR 10, VA Panel, 1ms MPRT, 165hz with Game Mode, Supports AMD FreeSync, HDMI and DisplayPort, 3 Years Local Warranty"
$ws.Range("B3").Value = "7,645"

$ws.Range("A4").Value = "Samsung 32-Inch QLED G8 Odyssey Gaming Monitor, with 1ms GtG Response time & 240Hz Refresh rate, Supports AMD FreeSync Premium Pro, Local Warranty"
$ws.Range("B4").Value = "82,678"

$ws.Range("A5").Value = "Samsung 27-Inch VA Gaming Monitor, 4ms GtG, 60hz with Game Mode, Supports AMD FreeSync, Local Warranty."
$ws.Range("B5").Value = "6,666"

$ws.Range("A6").Value = "Samsung 22-Inch IPS Gaming Monitor with Borderless Design, VGA and HDMI, 5ms GtG, 75hz with Game Mode, Supports AMD FreeSync, Local Warranty."
$ws.Range("B6").Value = "3,200"

# New rows 7-11, reusing the same item names with new prices
$ws.Range("A7").Value = 'Lenovo Monitor Legion R27fc-30 Gaming Curved Monitor, 27" FHD VA Display, Up to 280Hz Refresh Rate, 0.5ms (MPRT) Response Time, 1500R Curvature, FreeSync, Adaptive Sync & G-Sync Compatible, Black'
$ws.Range("B7").Value = "12,528"

$ws.Range("A8").Value = "Samsung 27-Inch G55C Odyssey QHD 2K Curved Gaming Monitor, HRDR 10, VA Panel, 1ms MPRT, 165hz with Game Mode, Supports AMD FreeSync, HDMI and DisplayPort, 3 Years Local Warranty"
$ws.Range("B8").Value = "12,555"

$ws.Range("A9").Value = "Samsung 32-Inch QLED G8 Odyssey Gaming Monitor, with 1ms GtG Response time & 240Hz Refresh rate, Supports AMD FreeSync Premium Pro, Local Warranty"
$ws.Range("B9").Value = "9,299"

$ws.Range("A10").Value = "Samsung 27-Inch VA Gaming Monitor, 4ms GtG, 60hz with Game Mode, Supports AMD FreeSync, Local Warranty."
$ws.Range("B10").Value = "11,999"

$ws.Range("A11").Value = "Samsung 22-Inch IPS Gaming Monitor with Borderless Design, VGA and HDMI, 5ms GtG, 75hz with Game Mode, Supports AMD FreeSync, Local Warranty."
$ws.Range("B11").Value = "7,645"

# Restore the default (General) number formatting now that the text values
# are locked in, so the cells keep using the workbook's base style.
$priceRange.ClearFormats()
